$p = $ppt.ActivePresentation
$s = $p.Slides.Item(32)
$sh = $s.Shapes.Item(2)
$sh.TextFrame.TextRange.Text = "Traceability Graph"
